$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: mark as registered
$ws.Range("C17").Value = "済"

# Row 24: mark as registered
$ws.Range("C24").Value = "済"

# Row 25: append journal citation to the source cell, then mark as registered
$ws.Range("B25").Value = "**神奈川県衛生研究所** <br> [Prevalence and Characteristics of _Salmonella_ and _Campylobacter_ in Retail Poultry Meat in Japan](https://www.jstage.jst.go.jp/article/yoken/70/3/70_JJID.2016.164/_pdf/-char/en) <br> (Japanese Journal of Infectious Diseases, 70, 239-247, 2017)"
$ws.Range("C25").Value = "済"
